$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text '2025-11-09 Sunday' '2025-11-10 Monday'
Replace-Text '41×88=3608' '31×29=899'
Replace-Text '44×45=1980' '35×68=2380'
Replace-Text '62×65=4030' '77×12=924'
Replace-Text '31×38=1178' '77×33=2541'
Replace-Text '28×87=2436' '25×12=300'
Replace-Text '57×60=3420' '96×52=4992'
Replace-Text '35×12=420' '15×51=765'
Replace-Text '76×22=1672' '85×53=4505'
Replace-Text '86×12=1032' '34×35=1190'
Replace-Text '91×35=3185' '74×29=2146'
Replace-Text '28×52=1456' '70×62=4340'
Replace-Text '27×32=864' '24×80=1920'
Replace-Text '48×75=3600' '54×51=2754'
Replace-Text '80×27=2160' '43×85=3655'
Replace-Text '12×68=816' '19×42=798'
Replace-Text '79×41=3239' '65×89=5785'
Replace-Text '44×84=3696' '65×34=2210'
Replace-Text '36×75=2700' '87×49=4263'
Replace-Text '66×56=3696' '46×54=2484'
Replace-Text '80×35=2800' '66×82=5412'
Replace-Text '12×75=900' '53×15=795'
Replace-Text '27×80=2160' '43×15=645'
Replace-Text '87×80=6960' '44×15=660'
Replace-Text '12×70=840' '20×66=1320'
Replace-Text '69×57=3933' '27×42=1134'
